$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A103").Value = "sandrine.grossetete-lalami@curie.fr"
$ws.Range("A104").Value = "cantini@bio.ens.psl.eu"

$ws.Range("F103").Select()
